# Auto-generated script applying updated Leve market-price data to the Titan_Profits workbook.
# Each row below corresponds to a single Leve row whose currentAveragePrice / LevePrice /
# LeveProfit figures were refreshed by the scheduled data-sync runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 126412.266
$ws.Range("I15").Value = 126412.266
$ws.Range("K15").Value = 379236.798
$ws.Range("M15").Value = -379067.798

$ws.Range("H17").Value = 3640.4814
$ws.Range("J17").Value = 3640.4814
$ws.Range("L17").Value = 10921.4442
$ws.Range("N17").Value = -11257.4442

$ws.Range("H40").Value = 2093.5715
$ws.Range("I40").Value = 2107.1
$ws.Range("K40").Value = 2107.1
$ws.Range("M40").Value = -1932.1

$ws.Range("H113").Value = 3802.75
$ws.Range("I113").Value = 1735
$ws.Range("K113").Value = 1735
$ws.Range("M113").Value = 1519

$ws.Range("H132").Value = 51605.863
$ws.Range("I132").Value = 57406.61
$ws.Range("J132").Value = 25502.5
$ws.Range("K132").Value = 172219.83
$ws.Range("L132").Value = 76507.5
$ws.Range("M132").Value = -169689.83
$ws.Range("N132").Value = -81567.5

$ws.Range("H133").Value = 42312.31
$ws.Range("J133").Value = 42312.31
$ws.Range("L133").Value = 42312.31
$ws.Range("N133").Value = -52432.31

$ws.Range("H135").Value = 1462.7858
$ws.Range("I135").Value = 1640.0869
$ws.Range("J135").Value = 647.2
$ws.Range("K135").Value = 14760.7821
$ws.Range("L135").Value = 5824.8
$ws.Range("M135").Value = -12225.7821
$ws.Range("N135").Value = -10894.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 72087.14
$ws.Range("I5").Value = 100727
$ws.Range("J5").Value = 487.5
$ws.Range("K5").Value = 100727
$ws.Range("L5").Value = 487.5
$ws.Range("M5").Value = -100615
$ws.Range("N5").Value = -711.5

$ws.Range("H32").Value = 21231.148
$ws.Range("I32").Value = 3708.3455
$ws.Range("K32").Value = 3708.3455
$ws.Range("M32").Value = -3421.3455

$ws.Range("H45").Value = 885.5
$ws.Range("I45").Value = 885.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 885.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -508.5
$ws.Range("N45").Value = $null

$ws.Range("H74").Value = 13001.272
$ws.Range("I74").Value = 1625
$ws.Range("J74").Value = 19502
$ws.Range("K74").Value = 1625
$ws.Range("L74").Value = 19502
$ws.Range("M74").Value = -751
$ws.Range("N74").Value = -21250

$ws.Range("H77").Value = 13001.272
$ws.Range("I77").Value = 1625
$ws.Range("J77").Value = 19502
$ws.Range("K77").Value = 8125
$ws.Range("L77").Value = 97510
$ws.Range("M77").Value = -3757
$ws.Range("N77").Value = -106246

$ws.Range("H102").Value = 1133.3334
$ws.Range("I102").Value = 960
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 960
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 662
$ws.Range("N102").Value = -5244

$ws.Range("H133").Value = 64115.25
$ws.Range("J133").Value = 64115.25
$ws.Range("L133").Value = 64115.25
$ws.Range("N133").Value = -69175.25

$ws.Range("H139").Value = 52170
$ws.Range("J139").Value = 52170
$ws.Range("L139").Value = 52170
$ws.Range("N139").Value = -62450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 72087.14
$ws.Range("I4").Value = 100727
$ws.Range("J4").Value = 487.5
$ws.Range("K4").Value = 100727
$ws.Range("L4").Value = 487.5
$ws.Range("M4").Value = -100612
$ws.Range("N4").Value = -717.5

$ws.Range("H20").Value = 1670.25
$ws.Range("I20").Value = 1414.375
$ws.Range("J20").Value = 1926.125
$ws.Range("K20").Value = 1414.375
$ws.Range("L20").Value = 1926.125
$ws.Range("M20").Value = -1167.375
$ws.Range("N20").Value = -2420.125

$ws.Range("H59").Value = 53300
$ws.Range("J59").Value = 53300
$ws.Range("L59").Value = 53300
$ws.Range("N59").Value = -54994

$ws.Range("H86").Value = 6536.7896
$ws.Range("I86").Value = 1336
$ws.Range("J86").Value = 13687.875
$ws.Range("K86").Value = 1336
$ws.Range("L86").Value = 13687.875
$ws.Range("M86").Value = -213
$ws.Range("N86").Value = -15933.875

$ws.Range("H89").Value = 6536.7896
$ws.Range("I89").Value = 1336
$ws.Range("J89").Value = 13687.875
$ws.Range("K89").Value = 6680
$ws.Range("L89").Value = 68439.375
$ws.Range("M89").Value = -1064
$ws.Range("N89").Value = -79671.375

$ws.Range("H99").Value = 2386.1428
$ws.Range("I99").Value = 2342.1667
$ws.Range("J99").Value = 2650
$ws.Range("K99").Value = 2342.1667
$ws.Range("L99").Value = 2650
$ws.Range("M99").Value = -844.1667000000002
$ws.Range("N99").Value = -5646

$ws.Range("H105").Value = 184805.62
$ws.Range("I105").Value = 2799.6758
$ws.Range("J105").Value = 558928.9399999999
$ws.Range("K105").Value = 2799.6758
$ws.Range("L105").Value = 558928.9399999999
$ws.Range("M105").Value = -1052.6758
$ws.Range("N105").Value = -562422.9399999999

$ws.Range("H132").Value = 45520
$ws.Range("J132").Value = 45520
$ws.Range("L132").Value = 45520
$ws.Range("N132").Value = -55640

$ws.Range("H135").Value = 42289.75
$ws.Range("J135").Value = 42289.75
$ws.Range("L135").Value = 42289.75
$ws.Range("N135").Value = -52429.75

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws.Range("H140").Value = 71110
$ws.Range("J140").Value = 71110
$ws.Range("L140").Value = 71110
$ws.Range("N140").Value = -81470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.11111
$ws.Range("J2").Value = 37.285713
$ws.Range("L2").Value = 37.285713
$ws.Range("N2").Value = -263.285713

$ws.Range("H80").Value = 26318568
$ws.Range("I80").Value = 2879
$ws.Range("J80").Value = 125002400
$ws.Range("K80").Value = 2879
$ws.Range("L80").Value = 125002400
$ws.Range("M80").Value = -1881
$ws.Range("N80").Value = -125004396

$ws.Range("H83").Value = 26318568
$ws.Range("I83").Value = 2879
$ws.Range("J83").Value = 125002400
$ws.Range("K83").Value = 14395
$ws.Range("L83").Value = 625012000
$ws.Range("M83").Value = -9403
$ws.Range("N83").Value = -625021984

$ws.Range("H97").Value = 1134.3334
$ws.Range("I97").Value = 1013.4167
$ws.Range("J97").Value = 1295.5555
$ws.Range("K97").Value = 1013.4167
$ws.Range("L97").Value = 1295.5555
$ws.Range("M97").Value = -517.4167
$ws.Range("N97").Value = -2287.5555

$ws.Range("H138").Value = 75685.8
$ws.Range("J138").Value = 75685.8
$ws.Range("L138").Value = 75685.8
$ws.Range("N138").Value = -85965.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3331.96
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 3414.95
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3414.95
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -3638.95

$ws.Range("H126").Value = 3331.96
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 3414.95
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 10244.85
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -15184.85

$ws.Range("H132").Value = 5193.2
$ws.Range("I132").Value = 3600
$ws.Range("J132").Value = 5591.5
$ws.Range("K132").Value = 10800
$ws.Range("L132").Value = 16774.5
$ws.Range("M132").Value = -8270
$ws.Range("N132").Value = -21834.5

$ws.Range("H136").Value = 5959.6
$ws.Range("I136").Value = 5959.6
$ws.Range("K136").Value = 17878.8
$ws.Range("M136").Value = -15328.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2872.0715
$ws.Range("I132").Value = 2686.476
$ws.Range("J132").Value = 3428.8572
$ws.Range("K132").Value = 8059.428
$ws.Range("L132").Value = 10286.5716
$ws.Range("M132").Value = -5529.428
$ws.Range("N132").Value = -15346.5716

$ws.Range("H136").Value = 2877.8572
$ws.Range("I136").Value = 1729
$ws.Range("J136").Value = 5750
$ws.Range("K136").Value = 5187
$ws.Range("L136").Value = 17250
$ws.Range("M136").Value = -2637
$ws.Range("N136").Value = -22350
